$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = [double]"1.194336911590932E-05"
$ws.Range("E2").Value = [double]"1.194336911590932E-05"
$ws.Range("D3").Value = [double]"0.5111863404899493"
$ws.Range("E3").Value = [double]"0.5111863404899493"
$ws.Range("D4").Value = [double]"0.6463990397357803"
$ws.Range("E4").Value = [double]"0.6463990397357803"
$ws.Range("D5").Value = [double]"0.2943143595054112"
$ws.Range("E5").Value = [double]"0.2943143595054112"
$ws.Range("D6").Value = [double]"0.2564953040126339"
$ws.Range("E6").Value = [double]"0.2564953040126339"
$ws.Range("D7").Value = [double]"0.9999420889760546"
$ws.Range("E7").Value = [double]"5.791102394536729E-05"
$ws.Range("D8").Value = [double]"0.6917293301616004"
$ws.Range("E8").Value = [double]"0.3082706698383996"
$ws.Range("D9").Value = [double]"0.5903995597952789"
$ws.Range("E9").Value = [double]"0.4096004402047211"
$ws.Range("C10").Value = $false
$ws.Range("D10").Value = [double]"0.09503836603061321"
$ws.Range("E10").Value = [double]"0.9049616339693868"
$ws.Range("C11").Value = $false
$ws.Range("D11").Value = [double]"0.0564600434686782"
$ws.Range("E11").Value = [double]"0.9435399565313218"
$ws.Range("F11").Value = [double]"0.8523609042167664"
$ws.Range("G11").Value = [double]"0.6"
$ws.Range("D12").Value = [double]"1.261254158108405E-07"
$ws.Range("E12").Value = [double]"1.261254158108405E-07"
$ws.Range("D13").Value = [double]"0.9539727933305705"
$ws.Range("E13").Value = [double]"0.9539727933305705"
$ws.Range("D14").Value = [double]"0.0003231820774379899"
$ws.Range("E14").Value = [double]"0.0003231820774379899"
$ws.Range("D15").Value = [double]"4.026375932031479E-11"
$ws.Range("E15").Value = [double]"4.026375932031479E-11"
$ws.Range("D16").Value = [double]"0.25242152854668"
$ws.Range("E16").Value = [double]"0.25242152854668"
$ws.Range("D17").Value = [double]"0.9999999815895072"
$ws.Range("E17").Value = [double]"1.841049279693863E-08"
$ws.Range("C18").Value = $false
$ws.Range("D18").Value = [double]"0.1627816639883918"
$ws.Range("E18").Value = [double]"0.8372183360116082"
$ws.Range("C19").Value = $false
$ws.Range("D19").Value = [double]"0.4096100523875683"
$ws.Range("E19").Value = [double]"0.5903899476124317"
$ws.Range("C20").Value = $false
$ws.Range("D20").Value = [double]"0.007875475801427967"
$ws.Range("E20").Value = [double]"0.992124524198572"
$ws.Range("C21").Value = $false
$ws.Range("D21").Value = [double]"0.002131725010589395"
$ws.Range("E21").Value = [double]"0.9978682749894106"
$ws.Range("F21").Value = [double]"1.707248330116272"
$ws.Range("G21").Value = [double]"0.5"
